$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("F5").Value = 12
$ws.Range("H5").Value = 12

# Row 8
$ws.Range("E8").Value = 33
$ws.Range("F8").Value = 12
$ws.Range("H8").Value = 12

# Row 9
$ws.Range("E9").Value = 16

# Row 12
$ws.Range("E12").Value = 21
$ws.Range("F12").Value = 6
$ws.Range("H12").Value = 6

# Row 14
$ws.Range("E14").Value = 33
$ws.Range("F14").Value = 16
$ws.Range("H14").Value = 16

# Row 15
$ws.Range("E15").Value = 76
$ws.Range("F15").Value = 38
$ws.Range("H15").Value = 38

# Row 16
$ws.Range("F16").Value = 78
$ws.Range("H16").Value = 78
